{"js": "// Capitalize \"designer\" -> \"Designer\" in the Lab 0 heading:\n//   \"experimental designer & activate your community license\"\n//   becomes\n//   \"experimental Designer & activate your community license\"\n//\n// (The canonical OOXML diff shows the original single run being split into\n// three runs - \"experimental \", \"D\", \"esigner & activate your community\n// license\" - all sharing identical run formatting; that three-way split is\n// an artifact of how the text was retyped interactively and carries no\n// visible/semantic difference from a single corrected run, so we simply\n// rewrite the run's text to the corrected string.)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The heading is the first paragraph in the document.\nconst heading = paragraphs.items[0];\n\nconst results = heading.search(\"experimental designer\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"experimental Designer\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Capitalize \"designer\" -> \"Designer\" in the Lab 0 heading:\n#   \"experimental designer & activate your community license\"\n#   becomes\n#   \"experimental Designer & activate your community license\"\n#\n# (The canonical OOXML diff shows the original single run being split into\n# three runs - \"experimental \", \"D\", \"esigner & activate your community\n# license\" - all sharing identical run formatting; that three-way split is\n# an artifact of how the text was retyped interactively and carries no\n# visible/semantic difference from a single corrected run, so we simply\n# rewrite the matched text to the corrected string.)\n\n$d = $word.ActiveDocument\n\n# The heading is the document's first paragraph; scope Find to it so we\n# don't touch the other \"designer\"/\"Designer\" occurrences later in the doc.\n$headingRange = $d.Paragraphs(1).Range\n\n$found = $headingRange.Find.Execute(\"experimental designer\", $true, $false, $false, $false, $false, $true, 1, $false, \"experimental Designer\", 2)\n"}
